# Applies the content edits described by the commit:
#  - DBD!G12: tighten "N: text" legend to "N:text" (drop space after colon)
#  - DBD!G14: reformat CdCode/PostDepCode legend
#  - DBD!G16: tighten "N.text" -> "N:text" for the two status codes that used a period
#  - DBD row 14 grows from 2 to 3 wrapped lines, so its row height increases
#  - Selection moves from G16 to G21

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# G12: RepayCode legend - remove the space that followed each "N:" prefix
$ws.Range("G12").Value = "共用代碼檔 RepayCode`n1:匯款轉帳`n2:銀行扣款`n3:員工扣薪`n4:支票`n5:特約金`n6:人事特約金`n7:定存特約`n8:劃撥存款"

# G14: PostDepCode legend - reworded/reformatted into three lines
$ws.Range("G14").Value = "CdCode.PostDepCode`nP:存簿`nG:劃撥"
# Row 14 is not a custom height in the source file; it auto-grows to fit the
# now 3-line wrapped label (2 lines -> 3 lines, 16.2pt per line)
$ws.Rows.Item(14).RowHeight = 48.6

# G16: Status legend - "2.取消授權" / "8.授權失敗" -> "2:取消授權" / "8:授權失敗"
$ws.Range("G16").Value = "空白:未授權`n0:授權成功    `n1:停止使用    `n2:取消授權     `n8:授權失敗`n9:已送出授權"

# Update the selection to match the saved view (G21)
$ws.Range("G21").Select()
